$d = $word.ActiveDocument

$d.Content.Find.Execute("PAGE 37 EXAMPLE PROGRAM", $true, $false, $false, $false, $false,
                         $true, 1, $false, "EnumDemo", 2)
